# Natmi following Dr Hou advice
# Update the Ptprc-Mrc1 LR-pair table: recompute stats for the existing
# Sending/Target cluster combinations and add the new "sCs" sending-cluster rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Ptprc"
$ws.Cells.Item(2, 3).Value2 = "Mrc1"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 205.313027
$ws.Cells.Item(2, 8).Value2 = 615.9390810000001
$ws.Cells.Item(2, 9).Value2 = 0.435242422384838
$ws.Cells.Item(2, 10).Value2 = 0.435242422384838
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 14.71647433333334
$ws.Cells.Item(2, 14).Value2 = 44.14942300000001
$ws.Cells.Item(2, 15).Value2 = 0.1054330184450109
$ws.Cells.Item(2, 16).Value2 = 0.1054330184450109
$ws.Cells.Item(2, 17).Value2 = 3021.483892144475
$ws.Cells.Item(2, 18).Value2 = 27193.35502930027
$ws.Cells.Item(2, 19).Value2 = 0.04588892234735186
$ws.Cells.Item(2, 20).Value2 = 0.04588892234735186

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Ptprc"
$ws.Cells.Item(3, 3).Value2 = "Mrc1"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 205.313027
$ws.Cells.Item(3, 8).Value2 = 615.9390810000001
$ws.Cells.Item(3, 9).Value2 = 0.435242422384838
$ws.Cells.Item(3, 10).Value2 = 0.435242422384838
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 0.07385333333333334
$ws.Cells.Item(3, 14).Value2 = 0.22156
$ws.Cells.Item(3, 15).Value2 = 0.0005291063388682706
$ws.Cells.Item(3, 16).Value2 = 0.0005291063388682705
$ws.Cells.Item(3, 17).Value2 = 15.16305142070667
$ws.Cells.Item(3, 18).Value2 = 136.46746278636
$ws.Cells.Item(3, 19).Value2 = 0.0002302895246281991
$ws.Cells.Item(3, 20).Value2 = 0.000230289524628199

# Row 4: ECs -> M2
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Ptprc"
$ws.Cells.Item(4, 3).Value2 = "Mrc1"
$ws.Cells.Item(4, 4).Value2 = "M2"
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 205.313027
$ws.Cells.Item(4, 8).Value2 = 615.9390810000001
$ws.Cells.Item(4, 9).Value2 = 0.435242422384838
$ws.Cells.Item(4, 10).Value2 = 0.435242422384838
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 124.584409
$ws.Cells.Item(4, 14).Value2 = 373.753227
$ws.Cells.Item(4, 15).Value2 = 0.8925582306290469
$ws.Cells.Item(4, 16).Value2 = 0.8925582306290469
$ws.Cells.Item(4, 17).Value2 = 25578.80212879605
$ws.Cells.Item(4, 18).Value2 = 230209.2191591644
$ws.Cells.Item(4, 19).Value2 = 0.3884792064185113
$ws.Cells.Item(4, 20).Value2 = 0.3884792064185113

# Row 5: ECs -> sCs
$ws.Cells.Item(5, 1).Value2 = "ECs"
$ws.Cells.Item(5, 2).Value2 = "Ptprc"
$ws.Cells.Item(5, 3).Value2 = "Mrc1"
$ws.Cells.Item(5, 4).Value2 = "sCs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 205.313027
$ws.Cells.Item(5, 8).Value2 = 615.9390810000001
$ws.Cells.Item(5, 9).Value2 = 0.435242422384838
$ws.Cells.Item(5, 10).Value2 = 0.435242422384838
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 0.2065306666666666
$ws.Cells.Item(5, 14).Value2 = 0.6195919999999999
$ws.Cells.Item(5, 15).Value2 = 0.001479644587073792
$ws.Cells.Item(5, 16).Value2 = 0.001479644587073792
$ws.Cells.Item(5, 17).Value2 = 42.40343634166133
$ws.Cells.Item(5, 18).Value2 = 381.630927074952
$ws.Cells.Item(5, 19).Value2 = 0.0006440040943466107
$ws.Cells.Item(5, 20).Value2 = 0.0006440040943466107

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Ptprc"
$ws.Cells.Item(6, 3).Value2 = "Mrc1"
$ws.Cells.Item(6, 4).Value2 = "ECs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 0.186821
$ws.Cells.Item(6, 8).Value2 = 0.5604629999999999
$ws.Cells.Item(6, 9).Value2 = 0.0003960412341120362
$ws.Cells.Item(6, 10).Value2 = 0.0003960412341120362
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 14.71647433333334
$ws.Cells.Item(6, 14).Value2 = 44.14942300000001
$ws.Cells.Item(6, 15).Value2 = 0.1054330184450109
$ws.Cells.Item(6, 16).Value2 = 0.1054330184450109
$ws.Cells.Item(6, 17).Value2 = 2.749346451427667
$ws.Cells.Item(6, 18).Value2 = 24.744118062849
$ws.Cells.Item(6, 19).Value2 = [double]"4.175582274111919e-05"
$ws.Cells.Item(6, 20).Value2 = [double]"4.175582274111919e-05"

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Ptprc"
$ws.Cells.Item(7, 3).Value2 = "Mrc1"
$ws.Cells.Item(7, 4).Value2 = "FAPs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 0.186821
$ws.Cells.Item(7, 8).Value2 = 0.5604629999999999
$ws.Cells.Item(7, 9).Value2 = 0.0003960412341120362
$ws.Cells.Item(7, 10).Value2 = 0.0003960412341120362
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 0.07385333333333334
$ws.Cells.Item(7, 14).Value2 = 0.22156
$ws.Cells.Item(7, 15).Value2 = 0.0005291063388682706
$ws.Cells.Item(7, 16).Value2 = 0.0005291063388682705
$ws.Cells.Item(7, 17).Value2 = 0.01379735358666667
$ws.Cells.Item(7, 18).Value2 = 0.12417618228
$ws.Cells.Item(7, 19).Value2 = [double]"2.095479274218911e-07"
$ws.Cells.Item(7, 20).Value2 = [double]"2.09547927421891e-07"

# Row 8: FAPs -> M2
$ws.Cells.Item(8, 1).Value2 = "FAPs"
$ws.Cells.Item(8, 2).Value2 = "Ptprc"
$ws.Cells.Item(8, 3).Value2 = "Mrc1"
$ws.Cells.Item(8, 4).Value2 = "M2"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 0.186821
$ws.Cells.Item(8, 8).Value2 = 0.5604629999999999
$ws.Cells.Item(8, 9).Value2 = 0.0003960412341120362
$ws.Cells.Item(8, 10).Value2 = 0.0003960412341120362
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 124.584409
$ws.Cells.Item(8, 14).Value2 = 373.753227
$ws.Cells.Item(8, 15).Value2 = 0.8925582306290469
$ws.Cells.Item(8, 16).Value2 = 0.8925582306290469
$ws.Cells.Item(8, 17).Value2 = 23.274983873789
$ws.Cells.Item(8, 18).Value2 = 209.474854864101
$ws.Cells.Item(8, 19).Value2 = 0.0003534898631751831
$ws.Cells.Item(8, 20).Value2 = 0.0003534898631751831

# Row 9: FAPs -> sCs
$ws.Cells.Item(9, 1).Value2 = "FAPs"
$ws.Cells.Item(9, 2).Value2 = "Ptprc"
$ws.Cells.Item(9, 3).Value2 = "Mrc1"
$ws.Cells.Item(9, 4).Value2 = "sCs"
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 0.186821
$ws.Cells.Item(9, 8).Value2 = 0.5604629999999999
$ws.Cells.Item(9, 9).Value2 = 0.0003960412341120362
$ws.Cells.Item(9, 10).Value2 = 0.0003960412341120362
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 0.2065306666666666
$ws.Cells.Item(9, 14).Value2 = 0.6195919999999999
$ws.Cells.Item(9, 15).Value2 = 0.001479644587073792
$ws.Cells.Item(9, 16).Value2 = 0.001479644587073792
$ws.Cells.Item(9, 17).Value2 = 0.03858426567733333
$ws.Cells.Item(9, 18).Value2 = 0.3472583910959999
$ws.Cells.Item(9, 19).Value2 = [double]"5.860002683118988e-07"
$ws.Cells.Item(9, 20).Value2 = [double]"5.860002683118988e-07"

# Row 10: M2 -> ECs
$ws.Cells.Item(10, 1).Value2 = "M2"
$ws.Cells.Item(10, 2).Value2 = "Ptprc"
$ws.Cells.Item(10, 3).Value2 = "Mrc1"
$ws.Cells.Item(10, 4).Value2 = "ECs"
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 266.1765593333333
$ws.Cells.Item(10, 8).Value2 = 798.529678
$ws.Cells.Item(10, 9).Value2 = 0.5642668278730386
$ws.Cells.Item(10, 10).Value2 = 0.5642668278730386
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 14.71647433333334
$ws.Cells.Item(10, 14).Value2 = 44.14942300000001
$ws.Cells.Item(10, 15).Value2 = 0.1054330184450109
$ws.Cells.Item(10, 16).Value2 = 0.1054330184450109
$ws.Cells.Item(10, 17).Value2 = 3917.180503563978
$ws.Cells.Item(10, 18).Value2 = 35254.6245320758
$ws.Cells.Item(10, 19).Value2 = 0.05949235487104589
$ws.Cells.Item(10, 20).Value2 = 0.05949235487104589

# Row 11: M2 -> FAPs
$ws.Cells.Item(11, 1).Value2 = "M2"
$ws.Cells.Item(11, 2).Value2 = "Ptprc"
$ws.Cells.Item(11, 3).Value2 = "Mrc1"
$ws.Cells.Item(11, 4).Value2 = "FAPs"
$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 6).Value2 = 1
$ws.Cells.Item(11, 7).Value2 = 266.1765593333333
$ws.Cells.Item(11, 8).Value2 = 798.529678
$ws.Cells.Item(11, 9).Value2 = 0.5642668278730386
$ws.Cells.Item(11, 10).Value2 = 0.5642668278730386
$ws.Cells.Item(11, 11).Value2 = 3
$ws.Cells.Item(11, 12).Value2 = 1
$ws.Cells.Item(11, 13).Value2 = 0.07385333333333334
$ws.Cells.Item(11, 14).Value2 = 0.22156
$ws.Cells.Item(11, 15).Value2 = 0.0005291063388682706
$ws.Cells.Item(11, 16).Value2 = 0.0005291063388682705
$ws.Cells.Item(11, 17).Value2 = 19.65802616196445
$ws.Cells.Item(11, 18).Value2 = 176.92223545768
$ws.Cells.Item(11, 19).Value2 = 0.0002985571554407161
$ws.Cells.Item(11, 20).Value2 = 0.000298557155440716

# Row 12: M2 -> M2
$ws.Cells.Item(12, 1).Value2 = "M2"
$ws.Cells.Item(12, 2).Value2 = "Ptprc"
$ws.Cells.Item(12, 3).Value2 = "Mrc1"
$ws.Cells.Item(12, 4).Value2 = "M2"
$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 6).Value2 = 1
$ws.Cells.Item(12, 7).Value2 = 266.1765593333333
$ws.Cells.Item(12, 8).Value2 = 798.529678
$ws.Cells.Item(12, 9).Value2 = 0.5642668278730386
$ws.Cells.Item(12, 10).Value2 = 0.5642668278730386
$ws.Cells.Item(12, 11).Value2 = 3
$ws.Cells.Item(12, 12).Value2 = 1
$ws.Cells.Item(12, 13).Value2 = 124.584409
$ws.Cells.Item(12, 14).Value2 = 373.753227
$ws.Cells.Item(12, 15).Value2 = 0.8925582306290469
$ws.Cells.Item(12, 16).Value2 = 0.8925582306290469
$ws.Cells.Item(12, 17).Value2 = 33161.44933419677
$ws.Cells.Item(12, 18).Value2 = 298453.0440077709
$ws.Cells.Item(12, 19).Value2 = 0.5036410014890244
$ws.Cells.Item(12, 20).Value2 = 0.5036410014890244

# Row 13: M2 -> sCs
$ws.Cells.Item(13, 1).Value2 = "M2"
$ws.Cells.Item(13, 2).Value2 = "Ptprc"
$ws.Cells.Item(13, 3).Value2 = "Mrc1"
$ws.Cells.Item(13, 4).Value2 = "sCs"
$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 6).Value2 = 1
$ws.Cells.Item(13, 7).Value2 = 266.1765593333333
$ws.Cells.Item(13, 8).Value2 = 798.529678
$ws.Cells.Item(13, 9).Value2 = 0.5642668278730386
$ws.Cells.Item(13, 10).Value2 = 0.5642668278730386
$ws.Cells.Item(13, 11).Value2 = 3
$ws.Cells.Item(13, 12).Value2 = 1
$ws.Cells.Item(13, 13).Value2 = 0.2065306666666666
$ws.Cells.Item(13, 14).Value2 = 0.6195919999999999
$ws.Cells.Item(13, 15).Value2 = 0.001479644587073792
$ws.Cells.Item(13, 16).Value2 = 0.001479644587073792
$ws.Cells.Item(13, 17).Value2 = 54.97362225015288
$ws.Cells.Item(13, 18).Value2 = 494.7626002513759
$ws.Cells.Item(13, 19).Value2 = 0.0008349143575276408
$ws.Cells.Item(13, 20).Value2 = 0.0008349143575276408

# Row 14: sCs -> ECs
$ws.Cells.Item(14, 1).Value2 = "sCs"
$ws.Cells.Item(14, 2).Value2 = "Ptprc"
$ws.Cells.Item(14, 3).Value2 = "Mrc1"
$ws.Cells.Item(14, 4).Value2 = "ECs"
$ws.Cells.Item(14, 5).Value2 = 2
$ws.Cells.Item(14, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(14, 7).Value2 = 0.044676
$ws.Cells.Item(14, 8).Value2 = 0.134028
$ws.Cells.Item(14, 9).Value2 = [double]"9.470850801135487e-05"
$ws.Cells.Item(14, 10).Value2 = [double]"9.470850801135488e-05"
$ws.Cells.Item(14, 11).Value2 = 3
$ws.Cells.Item(14, 12).Value2 = 1
$ws.Cells.Item(14, 13).Value2 = 14.71647433333334
$ws.Cells.Item(14, 14).Value2 = 44.14942300000001
$ws.Cells.Item(14, 15).Value2 = 0.1054330184450109
$ws.Cells.Item(14, 16).Value2 = 0.1054330184450109
$ws.Cells.Item(14, 17).Value2 = 0.6574732073160001
$ws.Cells.Item(14, 18).Value2 = 5.917258865844001
$ws.Cells.Item(14, 19).Value2 = [double]"9.985403872060642e-06"
$ws.Cells.Item(14, 20).Value2 = [double]"9.985403872060644e-06"

# Row 15: sCs -> FAPs
$ws.Cells.Item(15, 1).Value2 = "sCs"
$ws.Cells.Item(15, 2).Value2 = "Ptprc"
$ws.Cells.Item(15, 3).Value2 = "Mrc1"
$ws.Cells.Item(15, 4).Value2 = "FAPs"
$ws.Cells.Item(15, 5).Value2 = 2
$ws.Cells.Item(15, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(15, 7).Value2 = 0.044676
$ws.Cells.Item(15, 8).Value2 = 0.134028
$ws.Cells.Item(15, 9).Value2 = [double]"9.470850801135487e-05"
$ws.Cells.Item(15, 10).Value2 = [double]"9.470850801135488e-05"
$ws.Cells.Item(15, 11).Value2 = 3
$ws.Cells.Item(15, 12).Value2 = 1
$ws.Cells.Item(15, 13).Value2 = 0.07385333333333334
$ws.Cells.Item(15, 14).Value2 = 0.22156
$ws.Cells.Item(15, 15).Value2 = 0.0005291063388682706
$ws.Cells.Item(15, 16).Value2 = 0.0005291063388682705
$ws.Cells.Item(15, 17).Value2 = 0.003299471520000001
$ws.Cells.Item(15, 18).Value2 = 0.02969524368
$ws.Cells.Item(15, 19).Value2 = [double]"5.011087193356424e-08"
$ws.Cells.Item(15, 20).Value2 = [double]"5.011087193356424e-08"

# Row 16: sCs -> M2
$ws.Cells.Item(16, 1).Value2 = "sCs"
$ws.Cells.Item(16, 2).Value2 = "Ptprc"
$ws.Cells.Item(16, 3).Value2 = "Mrc1"
$ws.Cells.Item(16, 4).Value2 = "M2"
$ws.Cells.Item(16, 5).Value2 = 2
$ws.Cells.Item(16, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(16, 7).Value2 = 0.044676
$ws.Cells.Item(16, 8).Value2 = 0.134028
$ws.Cells.Item(16, 9).Value2 = [double]"9.470850801135487e-05"
$ws.Cells.Item(16, 10).Value2 = [double]"9.470850801135488e-05"
$ws.Cells.Item(16, 11).Value2 = 3
$ws.Cells.Item(16, 12).Value2 = 1
$ws.Cells.Item(16, 13).Value2 = 124.584409
$ws.Cells.Item(16, 14).Value2 = 373.753227
$ws.Cells.Item(16, 15).Value2 = 0.8925582306290469
$ws.Cells.Item(16, 16).Value2 = 0.8925582306290469
$ws.Cells.Item(16, 17).Value2 = 5.565933056484001
$ws.Cells.Item(16, 18).Value2 = 50.09339750835601
$ws.Cells.Item(16, 19).Value2 = [double]"8.453285833613182e-05"
$ws.Cells.Item(16, 20).Value2 = [double]"8.453285833613184e-05"

# Row 17: sCs -> sCs
$ws.Cells.Item(17, 1).Value2 = "sCs"
$ws.Cells.Item(17, 2).Value2 = "Ptprc"
$ws.Cells.Item(17, 3).Value2 = "Mrc1"
$ws.Cells.Item(17, 4).Value2 = "sCs"
$ws.Cells.Item(17, 5).Value2 = 2
$ws.Cells.Item(17, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(17, 7).Value2 = 0.044676
$ws.Cells.Item(17, 8).Value2 = 0.134028
$ws.Cells.Item(17, 9).Value2 = [double]"9.470850801135487e-05"
$ws.Cells.Item(17, 10).Value2 = [double]"9.470850801135488e-05"
$ws.Cells.Item(17, 11).Value2 = 3
$ws.Cells.Item(17, 12).Value2 = 1
$ws.Cells.Item(17, 13).Value2 = 0.2065306666666666
$ws.Cells.Item(17, 14).Value2 = 0.6195919999999999
$ws.Cells.Item(17, 15).Value2 = 0.001479644587073792
$ws.Cells.Item(17, 16).Value2 = 0.001479644587073792
$ws.Cells.Item(17, 17).Value2 = 0.009226964063999999
$ws.Cells.Item(17, 18).Value2 = 0.083042676576
$ws.Cells.Item(17, 19).Value2 = [double]"1.401349312288361e-07"
$ws.Cells.Item(17, 20).Value2 = [double]"1.401349312288361e-07"

